$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$tbl = $s.Shapes.Item(2).Table
$tbl.ApplyStyle("{06366723-33BB-4BD5-8D9F-F27E693FACD2}")
